$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.38681531086944
$ws.Range("C2").Value = 17.57355236682683
$ws.Range("D2").Value = 6.009292065293865
$ws.Range("E2").Value = 7.774514202376145
$ws.Range("G2").Value = 3.675116028684679
$ws.Range("I2").Value = 29.35107731766757
$ws.Range("M2").Value = 18.42280294353653
$ws.Range("N2").Value = 19.77824548082043
$ws.Range("B3").Value = 18.76123206061898
$ws.Range("C3").Value = 16.82908844989909
$ws.Range("D3").Value = 5.89723883164923
$ws.Range("E3").Value = 7.716581039622087
$ws.Range("G3").Value = 3.680583528159862
$ws.Range("I3").Value = 29.09808750766401
$ws.Range("M3").Value = 18.11792247215401
$ws.Range("N3").Value = 19.79473547196474
$ws.Range("B4").Value = 18.37390927265285
$ws.Range("C4").Value = 16.36159185971481
$ws.Range("D4").Value = 5.829519043069825
$ws.Range("E4").Value = 7.683062430559606
$ws.Range("G4").Value = 3.684105462121306
$ws.Range("I4").Value = 28.94948036332909
$ws.Range("M4").Value = 17.93468553182011
$ws.Range("N4").Value = 19.80672557746484
$ws.Range("B5").Value = 18.21554462688747
$ws.Range("C5").Value = 16.16879278750461
$ws.Range("D5").Value = 5.802232910995939
$ws.Range("E5").Value = 7.669928047211701
$ws.Range("G5").Value = 3.685582341283712
$ws.Range("I5").Value = 28.89065498565777
$ws.Range("M5").Value = 17.86110822910147
$ws.Range("N5").Value = 19.8120780505478
$ws.Range("B6").Value = 18.18922497010944
$ws.Range("C6").Value = 16.13665046681704
$ws.Range("D6").Value = 5.797721919279118
$ws.Range("E6").Value = 7.667779015077697
$ws.Range("G6").Value = 3.685830098257883
$ws.Range("I6").Value = 28.88099288494044
$ws.Range("M6").Value = 17.84895960152745
$ws.Range("N6").Value = 19.812994917848
$ws.Range("B7").Value = 18.37177526415463
$ws.Range("C7").Value = 16.35900052422708
$ws.Range("D7").Value = 5.829149746258763
$ws.Range("E7").Value = 7.682883160514328
$ws.Range("G7").Value = 3.684125210867066
$ws.Range("I7").Value = 28.94867995441296
$ws.Range("M7").Value = 17.93368868918121
$ws.Range("N7").Value = 19.80679587776061
$ws.Range("B8").Value = 19.17194347716924
$ws.Range("C8").Value = 17.31920073375996
$ws.Range("D8").Value = 5.970454742293945
$ws.Range("E8").Value = 7.754118292600665
$ws.Range("G8").Value = 3.676967128385445
$ws.Range("I8").Value = 29.26247541253379
$ws.Range("M8").Value = 18.31691699191426
$ws.Range("N8").Value = 19.78354254283886
$ws.Range("B9").Value = 20.70421659761326
$ws.Range("C9").Value = 19.10732295006772
$ws.Range("D9").Value = 6.254445693708989
$ws.Range("E9").Value = 7.909672615443547
$ws.Range("G9").Value = 3.664228905059383
$ws.Range("I9").Value = 29.92925659588488
$ws.Range("M9").Value = 19.09550872045144
$ws.Range("N9").Value = 19.75286004494656
$ws.Range("B10").Value = 21.79370749178197
$ws.Range("C10").Value = 20.34908810029061
$ws.Range("D10").Value = 6.46506195955012
$ws.Range("E10").Value = 8.032999933275587
$ws.Range("G10").Value = 3.655648695779345
$ws.Range("I10").Value = 30.44765154136955
$ws.Range("M10").Value = 19.67821042148061
$ws.Range("N10").Value = 19.73957483297518
$ws.Range("B11").Value = 22.27896660747294
$ws.Range("C11").Value = 20.89600677000592
$ws.Range("D11").Value = 6.560847685001749
$ws.Range("E11").Value = 8.090917757679367
$ws.Range("G11").Value = 3.651911546222876
$ws.Range("I11").Value = 30.6890157414335
$ws.Range("M11").Value = 19.94439478828297
$ws.Range("N11").Value = 19.73557619609658
$ws.Range("B12").Value = 22.46105170155342
$ws.Range("C12").Value = 21.10036672648758
$ws.Range("D12").Value = 6.597079405435405
$ws.Range("E12").Value = 8.113097168863755
$ws.Range("G12").Value = 3.650520036515521
$ws.Range("I12").Value = 30.78115291472636
$ws.Range("M12").Value = 20.04525296121207
$ws.Range("N12").Value = 19.73435889182163
$ws.Range("B13").Value = 22.42191348319073
$ws.Range("C13").Value = 21.05647862863066
$ws.Range("D13").Value = 6.589278613277163
$ws.Range("E13").Value = 8.108309668566491
$ws.Range("G13").Value = 3.650818673708003
$ws.Range("I13").Value = 30.7612776719067
$ws.Range("M13").Value = 20.02353023569326
$ws.Range("N13").Value = 19.73460781043805
$ws.Range("B14").Value = 22.29398134431785
$ws.Range("C14").Value = 20.91287547177729
$ws.Range("D14").Value = 6.56382950155361
$ws.Range("E14").Value = 8.092737598516747
$ws.Range("G14").Value = 3.651796592736732
$ws.Range("I14").Value = 30.69658148062296
$ws.Range("M14").Value = 19.95269180909433
$ws.Range("N14").Value = 19.73547008132552
$ws.Range("B15").Value = 22.21539638219211
$ws.Range("C15").Value = 20.82455233263936
$ws.Range("D15").Value = 6.548234882244783
$ws.Range("E15").Value = 8.083231026237923
$ws.Range("G15").Value = 3.652398672445664
$ws.Range("I15").Value = 30.65704749181032
$ws.Range("M15").Value = 19.90930612897764
$ws.Range("N15").Value = 19.73603699554767
$ws.Range("B16").Value = 21.76177019764847
$ws.Range("C16").Value = 20.3129696741229
$ws.Range("D16").Value = 6.458798465719916
$ws.Range("E16").Value = 8.029250226046113
$ws.Range("G16").Value = 3.655896250553015
$ws.Range("I16").Value = 30.43198433213604
$ws.Range("M16").Value = 19.66082839045217
$ws.Range("N16").Value = 19.73987754449372
$ws.Range("B17").Value = 21.4806964205309
$ws.Range("C17").Value = 19.99440493178546
$ws.Range("D17").Value = 6.403899176554076
$ws.Range("E17").Value = 7.996590057471546
$ws.Range("G17").Value = 3.658084278646809
$ws.Range("I17").Value = 30.29529345889771
$ws.Range("M17").Value = 19.50860777081048
$ws.Range("N17").Value = 19.74275936231749
$ws.Range("B18").Value = 21.3180651234432
$ws.Range("C18").Value = 19.80949188416469
$ws.Range("D18").Value = 6.37232233565745
$ws.Range("E18").Value = 7.977976554493699
$ws.Range("G18").Value = 3.659358414577537
$ws.Range("I18").Value = 30.21719794634758
$ws.Range("M18").Value = 19.42116480068425
$ws.Range("N18").Value = 19.74460921141118
$ws.Range("B19").Value = 21.26284110315372
$ws.Range("C19").Value = 19.74659980547225
$ws.Range("D19").Value = 6.361632066755702
$ws.Range("E19").Value = 7.971704257282132
$ws.Range("G19").Value = 3.659792507600479
$ws.Range("I19").Value = 30.19084818253043
$ws.Range("M19").Value = 19.39158015059672
$ws.Range("N19").Value = 19.74526848004957
$ws.Range("B20").Value = 21.51071847522681
$ws.Range("C20").Value = 20.02849213659632
$ws.Range("D20").Value = 6.40974359366082
$ws.Range("E20").Value = 8.000049110459113
$ws.Range("G20").Value = 3.657849742135402
$ws.Range("I20").Value = 30.3097904782348
$ws.Range("M20").Value = 19.52480123263211
$ws.Range("N20").Value = 19.74243266458898
$ws.Range("B21").Value = 22.33160485160891
$ws.Range("C21").Value = 20.95513090951601
$ws.Range("D21").Value = 6.571305899374334
$ws.Range("E21").Value = 8.097304897548693
$ws.Range("G21").Value = 3.651508713533365
$ws.Range("I21").Value = 30.7155647766928
$ws.Range("M21").Value = 19.97349792381265
$ws.Range("N21").Value = 19.73520873124717
$ws.Range("B22").Value = 22.85828476708469
$ws.Range("C22").Value = 21.54467663743613
$ws.Range("D22").Value = 6.676645754868922
$ws.Range("E22").Value = 8.162301206139658
$ws.Range("G22").Value = 3.647502341533898
$ws.Range("I22").Value = 30.98503375835221
$ws.Range("M22").Value = 20.26705271809106
$ws.Range("N22").Value = 19.73221930529108
$ws.Range("B23").Value = 22.57813937074247
$ws.Range("C23").Value = 21.23154265498241
$ws.Range("D23").Value = 6.620458496524555
$ws.Range("E23").Value = 8.127484952488244
$ws.Range("G23").Value = 3.649628072484258
$ws.Range("I23").Value = 30.84084213976264
$ws.Range("M23").Value = 20.1103810908647
$ws.Range("N23").Value = 19.7336554178948
$ws.Range("B24").Value = 21.49714872409178
$ws.Range("C24").Value = 20.01308679918917
$ws.Range("D24").Value = 6.407101375073526
$ws.Range("E24").Value = 7.998484763027849
$ws.Range("G24").Value = 3.657955725587743
$ws.Range("I24").Value = 30.30323484540191
$ws.Range("M24").Value = 19.51747994074832
$ws.Range("N24").Value = 19.74257976352843
$ws.Range("B25").Value = 20.29514688180665
$ws.Range("C25").Value = 18.63532756477963
$ws.Range("D25").Value = 6.177114105661265
$ws.Range("E25").Value = 7.865954436815235
$ws.Range("G25").Value = 3.667537252770455
$ws.Range("I25").Value = 29.74365115702656
$ws.Range("M25").Value = 18.88258495423499
$ws.Range("N25").Value = 19.75954658430596
